$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray leftover value in H3 (orphan "cocksmal" shared-string
# reference from before) first so that the shared-strings table slot it
# occupied can be reclaimed/reused by the new strings we add below.
$ws.Range("H3").Value = ""

# New column headers: contactnumber / diagnoses / treatment
$ws.Range("G1").Value = "contactnumber"
$ws.Range("H1").Value = "diagnoses"
$ws.Range("I1").Value = "treatment"

# New patient data columns
$ws.Range("G2").Value = 91111111
$ws.Range("H2").Value = "NULL"
$ws.Range("I2").Value = "NULL"

$ws.Range("G3").Value = 91111111
$ws.Range("H3").Value = "NULL"
$ws.Range("I3").Value = "NULL"

$ws.Range("G4").Value = 91111111
$ws.Range("H4").Value = "NULL"
$ws.Range("I4").Value = "NULL"

# Update the active selection to match the saved workbook state
$ws.Range("P11").Select()
